$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = "b554cb5"
$ws.Range("A3").Value = "a24a28d"
$ws.Range("A4").Value = "3b458f7"
$ws.Range("A5").Value = "1bea8b2"
$ws.Range("A6").Value = "535ca68"
$ws.Range("A7").Value = "8257ed5"
$ws.Range("A8").Value = "50bff96"
$ws.Range("A10").Value = "8d9ceb0"
$ws.Range("A11").Value = "5884be1"
$ws.Range("A12").Value = "db7e18d"
$ws.Range("A13").Value = "0eccbea"
$ws.Range("A14").Value = "cac0364"
$ws.Range("A15").Value = "a411821"
$ws.Range("A16").Value = "bc3c93f"
$ws.Range("A17").Value = "782bdc2"
$ws.Range("A18").Value = "764abe1"
$ws.Range("A20").Value = "b8b4d88"
$ws.Range("A21").Value = "b635fb7"
$ws.Range("A22").Value = "fbdf7f0"
$ws.Range("A23").Value = "5fc7519"
$ws.Range("A24").Value = "fbc3c63"
$ws.Range("A25").Value = "44c3b59"
$ws.Range("A26").Value = "a81c662"
$ws.Range("A27").Value = "2cd2b92"
$ws.Range("A28").Value = "a5efbd7"
$ws.Range("A29").Value = "c297cc2"
$ws.Range("A30").Value = "f0aeb17"
$ws.Range("A31").Value = "00f89a0"
$ws.Range("A32").Value = "d08d124"
$ws.Range("A33").Value = "b6d380f"
$ws.Range("A34").Value = "e3dbb90"
$ws.Range("A35").Value = "aff06c9"
$ws.Range("A36").Value = "7540dcc"
$ws.Range("A37").Value = "a806d69"
$ws.Range("A38").Value = "015ec49"
$ws.Range("A40").Value = "e73a78d"
$ws.Range("A41").Value = "ead35c7"
$ws.Range("A42").Value = "897b4eb"
$ws.Range("A43").Value = "a679e19"
$ws.Range("A44").Value = "b935cb5"
$ws.Range("A45").Value = "7fdf5a3"
$ws.Range("A46").Value = "eb048e7"
$ws.Range("A47").Value = "3ad3244"
$ws.Range("A48").Value = "6f4f49b"
$ws.Range("A49").Value = "4e4b33c"
$ws.Range("A50").Value = "89114a0"
$ws.Range("A51").Value = "d29a4fa"

# Values that look like numbers need to be forced to text. Writing them
# directly would store a numeric cell, and forcing NumberFormat = "@" on the
# target cell creates a brand-new style record (changing its "s" index), so
# instead stage the text in a helper cell via a formula (guaranteeing a text
# result), then copy/paste-special *values* onto the target - this preserves
# the target cell's existing style untouched.
$helper = $ws.Range("ZZ1")
$helper.Formula = '="8881302"'
$helper.Copy()
$ws.Range("A9").PasteSpecial(-4163)
$helper.Formula = '="2298935"'
$helper.Copy()
$ws.Range("A19").PasteSpecial(-4163)
$helper.Formula = '="5615e04"'
$helper.Copy()
$ws.Range("A39").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = $false
